# Auto-generated update of leve-profit market data across all job sheets
# (values refreshed by the scheduled market-data runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19: Unbreak My Heart | Roof Tile
$ws.Range("H19").Value = 2050.25
$ws.Range("J19").Value = 2598.7
$ws.Range("L19").Value = 2598.7
$ws.Range("N19").Value = -2948.7

# Row 33: Glazed and Confused | Clear Glass Lens
$ws.Range("H33").Value = 1613.2188
$ws.Range("I33").Value = 1602.963
$ws.Range("K33").Value = 1602.963
$ws.Range("M33").Value = -1373.963

# Row 39: Riches'' Brew | Hi-Potion of Mind
$ws.Range("H39").Value = 270.76923
$ws.Range("I39").Value = 251.9375
$ws.Range("J39").Value = 300.9
$ws.Range("K39").Value = 755.8125
$ws.Range("L39").Value = 902.6999999999999
$ws.Range("M39").Value = -459.8125
$ws.Range("N39").Value = -1494.7

# Row 41: The Write Stuff | Enchanted Mythril Ink
$ws.Range("H41").Value = 640.6667
$ws.Range("I41").Value = 689.8182
$ws.Range("K41").Value = 689.8182
$ws.Range("M41").Value = -249.8182

# Row 70: Consecrating Congregation | Holy Water
$ws.Range("H70").Value = 3250
$ws.Range("I70").Value = 3250
$ws.Range("K70").Value = 9750
$ws.Range("M70").Value = -9480

# Row 73: Curbing the Contagion (L) | Holy Water
$ws.Range("H73").Value = 3250
$ws.Range("I73").Value = 3250
$ws.Range("K73").Value = 9750
$ws.Range("M73").Value = -8814

# Row 113: Amaro Kart | Starch Glue
$ws.Range("H113").Value = 6326.0435
$ws.Range("I113").Value = 6048.3335
$ws.Range("K113").Value = 6048.3335
$ws.Range("M113").Value = -2794.3335

# Row 116: Growing Up | Growth Formula Kappa
$ws.Range("H116").Value = 5204.3335
$ws.Range("I116").Value = 4587.7617
$ws.Range("K116").Value = 4587.7617
$ws.Range("M116").Value = -1145.7617

# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 1536.8889
$ws.Range("I132").Value = 1580.5294
$ws.Range("J132").Value = 795
$ws.Range("K132").Value = 4741.5882
$ws.Range("L132").Value = 2385
$ws.Range("M132").Value = -2211.5882
$ws.Range("N132").Value = -7445

# Row 138: All-night Crafting | Cunning Craftsman''s Tisane
$ws.Range("H138").Value = 3782.02
$ws.Range("J138").Value = 3791.8572
$ws.Range("L138").Value = 11375.5716
$ws.Range("N138").Value = -21655.5716

$ws = $wb.Worksheets.Item("ARM")
# Row 31: I Was a Teenage Wailer | Iron Alembic
$ws.Range("H31").Value = 5138.5713
$ws.Range("I31").Value = 5138.5713
$ws.Range("K31").Value = 5138.5713
$ws.Range("M31").Value = -4844.5713

# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 9707.68
$ws.Range("I32").Value = 4035.038
$ws.Range("K32").Value = 4035.038
$ws.Range("M32").Value = -3748.038

# Row 45: Hollow Hallmarks | Mythril Ingot
$ws.Range("H45").Value = 29413914
$ws.Range("I45").Value = 31252096
$ws.Range("J45").Value = 3014
$ws.Range("K45").Value = 31252096
$ws.Range("L45").Value = 3014
$ws.Range("M45").Value = -31251719
$ws.Range("N45").Value = -3768

# Row 97: Ore for Me | High Steel Ingot
$ws.Range("H97").Value = 1550
$ws.Range("I97").Value = 1550
$ws.Range("K97").Value = 1550
$ws.Range("M97").Value = -1054

# Row 112: Wrapped Knuckles | Deepgold Gloves of Fending
$ws.Range("H112").Value = 97299
$ws.Range("J112").Value = 97299
$ws.Range("L112").Value = 97299
$ws.Range("N112").Value = -100253

# Row 132: Don''t Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 4147.8604
$ws.Range("I132").Value = 2504.0857
$ws.Range("K132").Value = 7512.257100000001
$ws.Range("M132").Value = -4982.257100000001

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 2647914.5
$ws.Range("I134").Value = 1492.7273
$ws.Range("J134").Value = 14292170
$ws.Range("K134").Value = 4478.1819
$ws.Range("L134").Value = 42876510
$ws.Range("M134").Value = -1943.1819
$ws.Range("N134").Value = -42881580

$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania''s Got Talent | Maple Lumber
$ws.Range("H7").Value = 205.5
$ws.Range("J7").Value = 276
$ws.Range("L7").Value = 276
$ws.Range("N7").Value = -502

# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 939975.8
$ws.Range("I31").Value = 11793.3125
$ws.Range("K31").Value = 11793.3125
$ws.Range("M31").Value = -11498.3125

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 939975.8
$ws.Range("I34").Value = 11793.3125
$ws.Range("K34").Value = 11793.3125
$ws.Range("M34").Value = -11591.3125

# Row 105: Zelkova, My Love | Zelkova Lumber
$ws.Range("H105").Value = 845
$ws.Range("I105").Value = 762.5
$ws.Range("K105").Value = 762.5
$ws.Range("M105").Value = 984.5

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 5759.0454
$ws.Range("I132").Value = 2417.5625
$ws.Range("K132").Value = 7252.6875
$ws.Range("M132").Value = -4722.6875

# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 2155.4878
$ws.Range("I134").Value = 1083.1177
$ws.Range("K134").Value = 3249.3531
$ws.Range("M134").Value = -714.3531000000003

# Row 141: No Greater Treasure | Claro Walnut Necklace of Gathering
$ws.Range("H141").Value = 276315.75
$ws.Range("J141").Value = 300195.53
$ws.Range("L141").Value = 300195.53
$ws.Range("N141").Value = -310555.53

$ws = $wb.Worksheets.Item("CUL")
# Row 25: Flakes for Friends | Apple Tart
$ws.Range("H25").Value = 119.666664
$ws.Range("I25").Value = 119.666664
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 358.999992
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -189.999992
$ws.Range("N25").ClearContents()

# Row 30: Picnic Panic | Apple Tart
$ws.Range("H30").Value = 119.666664
$ws.Range("I30").Value = 119.666664
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 358.999992
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -256.999992
$ws.Range("N30").ClearContents()

# Row 50: Moving Up in the World | Rolanberry Cheese
$ws.Range("H50").Value = 275
$ws.Range("I50").Value = 232
$ws.Range("J50").Value = 318
$ws.Range("K50").Value = 696
$ws.Range("L50").Value = 954
$ws.Range("M50").Value = -215
$ws.Range("N50").Value = -1916

# Row 53: Rolanberry Fields Forever | Rolanberry Cheese
$ws.Range("H53").Value = 275
$ws.Range("I53").Value = 232
$ws.Range("J53").Value = 318
$ws.Range("K53").Value = 696
$ws.Range("L53").Value = 954
$ws.Range("M53").Value = -215
$ws.Range("N53").Value = -1916

# Row 81: It Goes Down Smoothly | Frozen Spirits
$ws.Range("H81").Value = 1500
$ws.Range("I81").Value = 1500
$ws.Range("K81").Value = 4500
$ws.Range("M81").Value = -3377

# Row 84: Quenching the Flame (L) | Frozen Spirits
$ws.Range("H84").Value = 1500
$ws.Range("I84").Value = 1500
$ws.Range("K84").Value = 13500
$ws.Range("M84").Value = -7884

# Row 113: Can''t Eat Just One | Night Vinegar
$ws.Range("H113").Value = 1203.8334
$ws.Range("I113").Value = 720
$ws.Range("J113").Value = 1549.4286
$ws.Range("K113").Value = 2160
$ws.Range("L113").Value = 4648.2858
$ws.Range("M113").Value = 10
$ws.Range("N113").Value = -8988.2858

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 18976.2
$ws.Range("J131").Value = 27300
$ws.Range("L131").Value = 81900
$ws.Range("N131").Value = -91980

# Row 132: More Mezcal | Cooking Mezcal
$ws.Range("H132").Value = 1750
$ws.Range("J132").Value = 1500
$ws.Range("L132").Value = 13500
$ws.Range("N132").Value = -18560

# Row 133: Friends Are Food | Boiled Alpaca Steak
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 105: Untucked | Palladium Tuck
$ws.Range("H105").Value = 41278.168
$ws.Range("J105").Value = 41278.168
$ws.Range("L105").Value = 41278.168
$ws.Range("N105").Value = -48266.168

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 2715.08
$ws.Range("I132").Value = 2686.925
$ws.Range("J132").Value = 2827.7
$ws.Range("K132").Value = 8060.775000000001
$ws.Range("L132").Value = 8483.099999999999
$ws.Range("M132").Value = -5530.775000000001
$ws.Range("N132").Value = -13543.1

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban | Leather
$ws.Range("H7").Value = 859499.25
$ws.Range("I7").Value = 18534.4
$ws.Range("J7").Value = 1385102.2
$ws.Range("K7").Value = 18534.4
$ws.Range("L7").Value = 1385102.2
$ws.Range("M7").Value = -18422.4
$ws.Range("N7").Value = -1385326.2

# Row 22: Skin off Their Backs | Aldgoat Leather
$ws.Range("H22").Value = 983.625
$ws.Range("I22").Value = 974.4
$ws.Range("K22").Value = 974.4
$ws.Range("M22").Value = -679.4

# Row 27: Fire and Hide | Aldgoat Leather
$ws.Range("H27").Value = 983.625
$ws.Range("I27").Value = 974.4
$ws.Range("K27").Value = 974.4
$ws.Range("M27").Value = -867.4

# Row 46: Supply Side Logic | Boar Leather
$ws.Range("H46").Value = 3613.3333
$ws.Range("I46").Value = 2050
$ws.Range("K46").Value = 2050
$ws.Range("M46").Value = -1862

# Row 61: Spelling Me Softly | Raptor Leather
$ws.Range("H61").Value = 4198.4
$ws.Range("I61").Value = 3998
$ws.Range("K61").Value = 3998
$ws.Range("M61").Value = -3796

# Row 68: You Could Say It''s a Moving Target | Wyvern Leather
$ws.Range("H68").Value = 3499.5
$ws.Range("I68").Value = 4999
$ws.Range("K68").Value = 4999
$ws.Range("M68").Value = -4250

# Row 71: They Call It Bloody Mary (L) | Wyvern Leather
$ws.Range("H71").Value = 3499.5
$ws.Range("I71").Value = 4999
$ws.Range("K71").Value = 24995
$ws.Range("M71").Value = -21251

# Row 100: Tiger in the Sack | Tiger Leather
$ws.Range("H100").Value = 3434.5715
$ws.Range("I100").Value = 2823
$ws.Range("K100").Value = 2823
$ws.Range("M100").Value = -2282

# Row 101: A Stitch in Time | Marid Leather Gloves of Healing
$ws.Range("H101").Value = 75786.664
$ws.Range("J101").Value = 75786.664
$ws.Range("L101").Value = 75786.664
$ws.Range("N101").Value = -82276.664

# Row 106: If the Shoe Fits | Gazelleskin Boots of Casting
$ws.Range("H106").Value = 22399.8
$ws.Range("J106").Value = 22399.8
$ws.Range("L106").Value = 22399.8
$ws.Range("N106").Value = -24923.8

# Row 113: Peace in Rest | Atrociraptor Leather
$ws.Range("H113").Value = 4198.4
$ws.Range("I113").Value = 3998
$ws.Range("K113").Value = 3998
$ws.Range("M113").Value = -1828

# Row 126: Battered Books | Saiga Leather
$ws.Range("H126").Value = 859499.25
$ws.Range("I126").Value = 18534.4
$ws.Range("J126").Value = 1385102.2
$ws.Range("K126").Value = 55603.2
$ws.Range("L126").Value = 4155306.6
$ws.Range("M126").Value = -53133.2
$ws.Range("N126").Value = -4160246.6

$ws = $wb.Worksheets.Item("WVR")
# Row 19: Dirt Cheap | Stablehand''s Hat
$ws.Range("H19").Value = 2575
$ws.Range("J19").Value = 2575
$ws.Range("L19").Value = 2575
$ws.Range("N19").Value = -2923

# Row 113: A Tender Table | Pixie Floss
$ws.Range("H113").Value = 573.75
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 573.75
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1721.25
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6061.25

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 3210215.2
$ws.Range("I132").Value = 4733.4546
$ws.Range("K132").Value = 14200.3638
$ws.Range("M132").Value = -11670.3638
